$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to hold the company-name field.
$ws.Columns("B:B").Insert()

# New column B header + sample value (text formatted, matching the
# existing "@" text-format style used by the other text columns).
$ws.Range("B1").Value = "公司名称(与录入数据库名称一致)"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B2").Value = "阿里"

# Widen the new column to fit the longer header text.
$ws.Columns("B:B").ColumnWidth = 30.25

# Move the active selection to match the edited workbook's saved view.
$ws.Range("H8").Select()
